$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 179308
$ws.Range("I6").Value = 253662.67
$ws.Range("K6").Value = 760988.01
$ws.Range("M6").Value = -760876.01
$ws.Range("H17").Value = 222.47826
$ws.Range("J17").Value = 222.47826
$ws.Range("L17").Value = 667.43478
$ws.Range("N17").Value = -1003.43478
$ws.Range("H100").Value = 1766.8
$ws.Range("I100").Value = 1977.5
$ws.Range("J100").Value = 1626.3334
$ws.Range("K100").Value = 1977.5
$ws.Range("L100").Value = 1626.3334
$ws.Range("M100").Value = -1436.5
$ws.Range("N100").Value = -2708.3334
$ws.Range("H113").Value = 53973.527
$ws.Range("I113").Value = 84553.164
$ws.Range("J113").Value = 1551.2858
$ws.Range("K113").Value = 84553.164
$ws.Range("L113").Value = 1551.2858
$ws.Range("M113").Value = -81299.164
$ws.Range("N113").Value = -8059.2858
$ws.Range("H116").Value = 6690.2856
$ws.Range("I116").Value = 8266.4
$ws.Range("J116").Value = 2750
$ws.Range("K116").Value = 8266.4
$ws.Range("L116").Value = 2750
$ws.Range("M116").Value = -4824.4
$ws.Range("N116").Value = -9634
$ws.Range("H132").Value = 6761716
$ws.Range("I132").Value = 9264934
$ws.Range("J132").Value = 3029
$ws.Range("K132").Value = 27794802
$ws.Range("L132").Value = 9087
$ws.Range("M132").Value = -27792272
$ws.Range("N132").Value = -14147

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1940.3667
$ws.Range("I74").Value = 1008.3077
$ws.Range("J74").Value = 7998.75
$ws.Range("K74").Value = 1008.3077
$ws.Range("L74").Value = 7998.75
$ws.Range("M74").Value = -134.3077
$ws.Range("N74").Value = -9746.75
$ws.Range("H77").Value = 1940.3667
$ws.Range("I77").Value = 1008.3077
$ws.Range("J77").Value = 7998.75
$ws.Range("K77").Value = 5041.5385
$ws.Range("L77").Value = 39993.75
$ws.Range("M77").Value = -673.5384999999997
$ws.Range("N77").Value = -48729.75
$ws.Range("H98").Value = 8333
$ws.Range("J98").Value = 8333
$ws.Range("L98").Value = 8333
$ws.Range("N98").Value = -14323
$ws.Range("H102").Value = 40804.08
$ws.Range("I102").Value = 73848.5
$ws.Range("K102").Value = 73848.5
$ws.Range("M102").Value = -72226.5
$ws.Range("H106").Value = 42980
$ws.Range("J106").Value = 42980
$ws.Range("L106").Value = 42980
$ws.Range("N106").Value = -45504

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2126.0557
$ws.Range("I134").Value = 2119.9395
$ws.Range("J134").Value = 2193.3333
$ws.Range("K134").Value = 6359.818499999999
$ws.Range("L134").Value = 6579.999899999999
$ws.Range("M134").Value = -3824.818499999999
$ws.Range("N134").Value = -11649.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1123.75
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1247.5
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1247.5
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -1821.5
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H31").Value = 22170.12
$ws.Range("I31").Value = 34423.133
$ws.Range("J31").Value = 3790.6
$ws.Range("K31").Value = 34423.133
$ws.Range("L31").Value = 3790.6
$ws.Range("M31").Value = -34128.133
$ws.Range("N31").Value = -4380.6
$ws.Range("H34").Value = 22170.12
$ws.Range("I34").Value = 34423.133
$ws.Range("J34").Value = 3790.6
$ws.Range("K34").Value = 34423.133
$ws.Range("L34").Value = 3790.6
$ws.Range("M34").Value = -34221.133
$ws.Range("N34").Value = -4194.6
$ws.Range("H113").Value = 1123.75
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1247.5
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1247.5
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5587.5
$ws.Range("H122").Value = 921.2
$ws.Range("I122").Value = 898
$ws.Range("J122").Value = 1014
$ws.Range("K122").Value = 2694
$ws.Range("L122").Value = 3042
$ws.Range("M122").Value = -244
$ws.Range("N122").Value = -7942
$ws.Range("H134").Value = 1542.7931
$ws.Range("I134").Value = 1537.591
$ws.Range("K134").Value = 4612.772999999999
$ws.Range("M134").Value = -2077.772999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 47.444443
$ws.Range("J12").Value = 53.25
$ws.Range("L12").Value = 159.75
$ws.Range("N12").Value = -505.75
$ws.Range("H23").Value = 616.1875
$ws.Range("I23").Value = 222.25
$ws.Range("J23").Value = 747.5
$ws.Range("K23").Value = 666.75
$ws.Range("L23").Value = 2242.5
$ws.Range("M23").Value = -431.75
$ws.Range("N23").Value = -2712.5
$ws.Range("H36").Value = 250.5
$ws.Range("I36").Value = 250.5
$ws.Range("K36").Value = 751.5
$ws.Range("M36").Value = -582.5
$ws.Range("H58").Value = 1285.5714
$ws.Range("I58").Value = 1005
$ws.Range("J58").Value = 1332.3334
$ws.Range("K58").Value = 3015
$ws.Range("L58").Value = 3997.0002
$ws.Range("M58").Value = -2887
$ws.Range("N58").Value = -4253.0002
$ws.Range("H76").Value = 2416.6667
$ws.Range("H79").Value = 2416.6667
$ws.Range("H131").Value = 744.63
$ws.Range("I131").Value = 450.33334
$ws.Range("J131").Value = 796.5647
$ws.Range("K131").Value = 1351.00002
$ws.Range("L131").Value = 2389.6941
$ws.Range("M131").Value = 3688.99998
$ws.Range("N131").Value = -12469.6941
$ws.Range("H132").Value = 5962.75
$ws.Range("J132").Value = 6040.5
$ws.Range("L132").Value = 54364.5
$ws.Range("N132").Value = -59424.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 195587.67
$ws.Range("I102").Value = 2094.5
$ws.Range("K102").Value = 2094.5
$ws.Range("M102").Value = -472.5
$ws.Range("H107").Value = 1122701.5
$ws.Range("I107").Value = 438.85715
$ws.Range("J107").Value = 5050621
$ws.Range("K107").Value = 438.85715
$ws.Range("L107").Value = 5050621
$ws.Range("M107").Value = 1481.14285
$ws.Range("N107").Value = -5054461
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170
$ws.Range("H122").Value = 751.75
$ws.Range("I122").Value = 703.5
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 2110.5
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = 339.5
$ws.Range("N122").Value = -7300
$ws.Range("H132").Value = 2127.5
$ws.Range("I132").Value = 1789.9259
$ws.Range("J132").Value = 3140.2222
$ws.Range("K132").Value = 5369.7777
$ws.Range("L132").Value = 9420.6666
$ws.Range("M132").Value = -2839.7777
$ws.Range("N132").Value = -14480.6666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3745.25
$ws.Range("I122").Value = 3680.6
$ws.Range("J122").Value = 3853
$ws.Range("K122").Value = 11041.8
$ws.Range("L122").Value = 11559
$ws.Range("M122").Value = -8591.799999999999
$ws.Range("N122").Value = -16459

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6286
$ws.Range("I54").Value = 4535
$ws.Range("J54").Value = 6452.7617
$ws.Range("K54").Value = 4535
$ws.Range("L54").Value = 6452.7617
$ws.Range("M54").Value = -4015
$ws.Range("N54").Value = -7492.7617
$ws.Range("H81").Value = 251081.38
$ws.Range("I81").Value = 333790
$ws.Range("J81").Value = 201456.2
$ws.Range("K81").Value = 667580
$ws.Range("L81").Value = 402912.4
$ws.Range("M81").Value = -666519
$ws.Range("N81").Value = -405034.4
$ws.Range("H84").Value = 251081.38
$ws.Range("I84").Value = 333790
$ws.Range("J84").Value = 201456.2
$ws.Range("K84").Value = 3337900
$ws.Range("L84").Value = 2014562
$ws.Range("M84").Value = -3332596
$ws.Range("N84").Value = -2025170
